$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new ifo GDP component analysis preprocessing values
$ws.Range("D88").Value = 0.791995474
$ws.Range("D89").Value = 0.788120887
$ws.Range("D90").Value = 0.597740902
$ws.Range("D91").Value = 0.620527487
$ws.Range("C92").Value = 0.241887844
$ws.Range("C93").Value = 0.331651578
$ws.Range("C94").Value = 0.154182215
$ws.Range("C95").Value = 0.166899468
$ws.Range("C96").Value = 0.042359665
$ws.Range("C97").Value = 0.266698307
